# Updates activity report figures (hourly stats refreshed, chart now starts from row 8 / real-time line)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - AMM
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 516.237
$ws.Range("D2").Value = 735.242
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 139
$ws.Range("I2").Value = 2432
$ws.Range("J2").Value = -69.76800986842106

# Row 3 - IPR
$ws.Range("C3").Value = 165
$ws.Range("D3").Value = 168
$ws.Range("I3").Value = 216
$ws.Range("J3").Value = -22.22222222222222

# Row 4 - MIG
$ws.Range("C4").Value = 290
$ws.Range("D4").Value = 296
$ws.Range("G4").Value = 9
$ws.Range("I4").Value = 510
$ws.Range("J4").Value = -41.96078431372548

# Row 5 - MOB
$ws.Range("C5").Value = 904
$ws.Range("D5").Value = 1026
$ws.Range("G5").Value = 18
$ws.Range("H5").Value = 74
$ws.Range("I5").Value = 1434
$ws.Range("J5").Value = -28.45188284518828

# Row 6 - MOB PRE
$ws.Range("C6").Value = 510
$ws.Range("D6").Value = 554
$ws.Range("E6").Value = 44
$ws.Range("I6").Value = 1232
$ws.Range("J6").Value = -55.03246753246754

# Row 7 - MSK
$ws.Range("C7").Value = 256
$ws.Range("D7").Value = 286
$ws.Range("I7").Value = 324
$ws.Range("J7").Value = -11.72839506172839

# Row 8 - TEC
$ws.Range("C8").Value = 437
$ws.Range("D8").Value = 441
$ws.Range("I8").Value = 1013
$ws.Range("J8").Value = -56.4659427443238

# Row 9 - TST
$ws.Range("C9").Value = 148
$ws.Range("D9").Value = 156
$ws.Range("G9").Value = 5
$ws.Range("I9").Value = 164
$ws.Range("J9").Value = -4.878048780487809

# Row 10 - VIP
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 8
$ws.Range("J10").Value = 100

# Row 11 - WLC
$ws.Range("C11").Value = 35
$ws.Range("D11").Value = 36
$ws.Range("G11").Value = 1
$ws.Range("I11").Value = 111
$ws.Range("J11").Value = -67.56756756756756
